$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 "Marking": Right points and Wrong points corrected
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -2

# Row 12 "Total": recompute totals & display string
$ws.Range("B12").Value = 36
$ws.Range("C12").Value = -12
$ws.Range("E12").Value = "24 / 112"
